# homepage header and middle content arrow design changes suggested by Mylesi
#
# Adds a new "Sr. No. 6" task row (row 11) to the "log sheet" worksheet:
#   - Task: "Homepage deisgn changes in header and middle content arrow and integration"
#   - Website: myguitarpal
#   - Date: 30/09/2013
#   - Time(Hrs): 3, Rate: 13, Cost: 39, Currency: USD
# Row 11 gets the same 30pt height used by the other multi-line task row (row 9).
# The totals (Total Hours Consumed / Total Used Payment) recalc automatically.
# The active cell/selection moves from G10 to G11 to reflect the newly entered row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("log sheet")

$ws.Range("A11").Value = 6
$ws.Range("B11").Value = "Homepage deisgn changes in header and middle content arrow and integration"
$ws.Range("C11").Value = "myguitarpal"
$ws.Range("D11").Value = "30/09/2013"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 13
$ws.Range("G11").Value = 39
$ws.Range("H11").Value = "USD"

$ws.Rows.Item(11).RowHeight = 30

$ws.Range("G11").Select()
